$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete "J" series row (row 10) and the "L" series row (row 12).
# The former "K" series row (row 11) shifts up to become the new row 10.
# Delete from the bottom up so row indices stay valid.
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(10).Delete()

# Update header label in K1
$ws.Range("K1").Value = "Design Strength"

# Row 2 ("B" series) - updated analysis results
$ws.Range("D2").Value = 133.465852319
$ws.Range("E2").Value = 37.12289028844057
$ws.Range("F2").Value = 27.81452307344672
$ws.Range("I2").Value = 147.6313386352384
$ws.Range("J2").Value = 3.732721538111613
$ws.Range("K2").Value = 43.04876000792099
$ws.Range("N2").Value = 137.5265212300033
$ws.Range("O2").Value = 300.4577935658554
$ws.Range("P2").Value = 680.7369728571302
$ws.Range("Q2").Value = 552.1467546469989

# Row 3 ("C" series) - updated analysis results
$ws.Range("D3").Value = 218.92485695
$ws.Range("E3").Value = 63.71010620895059
$ws.Range("F3").Value = 29.10135792551929
$ws.Range("I3").Value = 243.0968542786452
$ws.Range("J3").Value = 3.462820561380488
$ws.Range("K3").Value = 64.3939687516219
$ws.Range("N3").Value = 15.59242455167278
$ws.Range("O3").Value = 104.6741842613235
$ws.Range("P3").Value = 931.2633461952848
$ws.Range("Q3").Value = 537.3673108146872

# Row 4 ("D" series) - updated analysis results
$ws.Range("D4").Value = 126.70803296
$ws.Range("E4").Value = 18.0425110891256
$ws.Range("F4").Value = 14.23943744341875
$ws.Range("I4").Value = 134.1934515757539
$ws.Range("J4").Value = 7.827709892680903
$ws.Range("K4").Value = 74.56024459595899
$ws.Range("N4").Value = 28.77239337819473
$ws.Range("O4").Value = 142.7393917410353
$ws.Range("P4").Value = 252.9164449060363
$ws.Range("Q4").Value = 206.1186936916146

# Row 5 ("E" series) - updated analysis results
$ws.Range("D5").Value = 77.15356611749999
$ws.Range("E5").Value = 9.069186090205813
$ws.Range("F5").Value = 11.75472054835938
$ws.Range("I5").Value = 81.08861249799605
$ws.Range("J5").Value = 8.678824106185585
$ws.Range("K5").Value = 47.72711205448898
$ws.Range("N5").Value = 276.8699166878031
$ws.Range("O5").Value = 366.6020013147311
$ws.Range("P5").Value = 160.1083736445783
$ws.Range("Q5").Value = 155.0123028439635

# Row 6 ("F" series) - updated analysis results
$ws.Range("D6").Value = 59.56530107250001
$ws.Range("E6").Value = 12.80282291835289
$ws.Range("F6").Value = 21.49376010501468
$ws.Range("I6").Value = 65.07948742504992
$ws.Range("J6").Value = 4.431981796208421
$ws.Range("K6").Value = 23.04998999891203
$ws.Range("N6").Value = 1712.765360158157
$ws.Range("O6").Value = 909.5129508948302
$ws.Range("P6").Value = 302.7249032640486
$ws.Range("Q6").Value = 349.1982764246479

# Row 7 ("G" series) - updated fracture_stress_mean only
$ws.Range("D7").Value = 69.21942002999999

# Row 8 ("H" series) - updated fracture_stress_mean only
$ws.Range("D8").Value = 74.62428989

# Row 9 ("I" series) - updated fracture_stress_mean only
$ws.Range("D9").Value = 66.09967249

# Row 10 ("K" series, formerly row 11) - updated fracture_stress_mean only
$ws.Range("D10").Value = 89.27432161999999
